$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 278
$ws1.Range("F4").Value = 197
$ws1.Range("F5").Value = 838
$ws1.Range("F7").Value = 535
$ws1.Range("F8").Value = 102
$ws1.Range("F9").Value = 548
$ws1.Range("F10").Value = 504
$ws1.Range("F12").Value = 36
$ws1.Range("F14").Value = 212

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 112
$ws2.Range("F9").Value = 187

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 6250
$ws3.Range("F4").Value = 1859

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 6250
$ws4.Range("F4").Value = 1859
$ws4.Range("F5").Value = 278
$ws4.Range("F10").Value = 197
$ws4.Range("F13").Value = 838
$ws4.Range("F14").Value = 112
$ws4.Range("F17").Value = 535
$ws4.Range("F18").Value = 187
$ws4.Range("F19").Value = 102
$ws4.Range("F20").Value = 548
$ws4.Range("F22").Value = 504
$ws4.Range("F26").Value = 36
$ws4.Range("F33").Value = 212
